$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 6 new daily data rows (124-129, 2024-09-19 .. 2024-09-24) below the
# existing data, extending the sheet's used range from A1:Z123 to A1:Z129.
# (Scientific-notation literals like 1.9056E-06 aren't tokenized by this PS
# parser, so those few values are routed through a [double] string cast.)

$ws.Range("A124").Value = 45554
$ws.Range("B124").Value = 557.1388109322
$ws.Range("C124").Value = 170.9286633835
$ws.Range("I124").Value = 255.836414559
$ws.Range("K124").Value = 294.135484908345
$ws.Range("N124").Value = 41.53587124736
$ws.Range("O124").Value = 0.02274804
$ws.Range("Q124").Value = [double]"1.9056E-06"
$ws.Range("U124").Value = 270.0581572768631
$ws.Range("Z124").Value = 314.231466513482
$ws.Range("A125").Value = 45555
$ws.Range("B125").Value = 559.3785893190001
$ws.Range("C125").Value = 177.59812689
$ws.Range("I125").Value = 262.586385046
$ws.Range("K125").Value = 295.10622908296
$ws.Range("N125").Value = 41.17499344224
$ws.Range("O125").Value = 0.02282828
$ws.Range("Q125").Value = [double]"1.98E-06"
$ws.Range("U125").Value = 269.2905831680705
$ws.Range("Z125").Value = 326.876030577782
$ws.Range("A126").Value = 45556
$ws.Range("B126").Value = 560.6877081888
$ws.Range("C126").Value = 181.13428074
$ws.Range("I126").Value = 267.510103571
$ws.Range("K126").Value = 295.10622908296
$ws.Range("N126").Value = 44.05037466368
$ws.Range("O126").Value = 0.023614632
$ws.Range("Q126").Value = [double]"2.0208E-06"
$ws.Range("U126").Value = 272.488808621373
$ws.Range("Z126").Value = 330.107419171992
$ws.Range("A127").Value = 45557
$ws.Range("B127").Value = 562.7216174328
$ws.Range("C127").Value = 178.95711935
$ws.Range("I127").Value = 258.951785553
$ws.Range("K127").Value = 294.523782578191
$ws.Range("N127").Value = 41.86182539392
$ws.Range("O127").Value = 0.023618644
$ws.Range("Q127").Value = [double]"1.9224E-06"
$ws.Range("U127").Value = 266.4761447691643
$ws.Range("Z127").Value = 321.087630139458
$ws.Range("A128").Value = 45558
$ws.Range("B128").Value = 560.6083166922
$ws.Range("C128").Value = 183.5312383595
$ws.Range("I128").Value = 259.238256449
$ws.Range("K128").Value = 295.688675587729
$ws.Range("N128").Value = 43.67785563904
$ws.Range("O128").Value = 0.02425254
$ws.Range("Q128").Value = [double]"1.9512E-06"
$ws.Range("U128").Value = 286.4330715977719
$ws.Range("Z128").Value = 329.910725953214
$ws.Range("A129").Value = 45559
$ws.Range("B129").Value = 568.775019906
$ws.Range("C129").Value = 183.96320382
$ws.Range("I129").Value = 273.454374663
$ws.Range("K129").Value = 293.747187238499
$ws.Range("N129").Value = 43.16564198016
$ws.Range("O129").Value = 0.024308708
$ws.Range("Q129").Value = [double]"2.0616E-06"
$ws.Range("U129").Value = 294.4925997400941
$ws.Range("Z129").Value = 324.993395483764

# Give the new dates in column A the same style as the existing date column
# (style index 2: bold/centered/bordered, numFmt "YYYY-MM-DD HH:MM:SS") by
# copying formats only (xlPasteFormats = -4122) from A123, so it reuses the
# existing cellXfs entry instead of minting a new style.
$ws.Range("A123").Copy()
$ws.Range("A124:A129").PasteSpecial(-4122)
$excel.CutCopyMode = $false
